{"js": "// Resume update: change the self-rated skill level from \"Novice\"/\"Beginner\"\n// to \"Basic\" for six entries in the \"Programming Languages/Frameworks\" list\n// (Python 3, C++, Bootstrap, Linux bash, JavaScript, Ruby on Rails), and\n// carry the \"_GoBack\" bookmark (Word's \"last edit position\" marker) to the\n// final edit made \u2014 right after \"Ruby on Rails (Basic\" \u2014 instead of its\n// original spot at the very start of the document.\n\n// 1. \"Novice\" -> \"Basic\" (Python 3, C++, Bootstrap)\nconst noviceResults = context.document.body.search(\"Novice\", { matchCase: true });\nnoviceResults.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < noviceResults.items.length; i++) {\n  noviceResults.items[i].insertText(\"Basic\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2. \"Beginner\" -> \"Basic\" (Linux bash, JavaScript, Ruby on Rails)\nconst beginnerResults = context.document.body.search(\"Beginner\", { matchCase: true });\nbeginnerResults.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < beginnerResults.items.length; i++) {\n  beginnerResults.items[i].insertText(\"Basic\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 3. Move the \"_GoBack\" bookmark (present at the top of the resume) to sit\n//    right after \"Ruby on Rails (Basic\" \u2014 the location of the last text\n//    edit \u2014 which is where Word leaves it after this kind of change.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst rubyLineResults = context.document.body.search(\"Ruby on Rails (Basic)\", { matchCase: true });\nrubyLineResults.load(\"text\");\nawait context.sync();\nif (rubyLineResults.items.length > 0) {\n  const basicInLine = rubyLineResults.items[0].search(\"Basic\", { matchCase: true });\n  basicInLine.load(\"text\");\n  await context.sync();\n  if (basicInLine.items.length > 0) {\n    basicInLine.items[0].getRange(\"End\").insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n", "ps1": "# Resume update: change the self-rated skill level from \"Novice\"/\"Beginner\"\n# to \"Basic\" for six entries in the \"Programming Languages/Frameworks\" list\n# (Python 3, C++, Bootstrap, Linux bash, JavaScript, Ruby on Rails), and\n# carry the \"_GoBack\" bookmark (Word's \"last edit position\" marker) to the\n# final edit made -- right after \"Ruby on Rails (Basic\" -- instead of its\n# original spot at the very start of the document.\n\n$d = $word.ActiveDocument\n\n# 1. \"Novice\" -> \"Basic\" (Python 3, C++, Bootstrap)\n$findNovice = $d.Content.Find\n$findNovice.Text = \"Novice\"\n$findNovice.Replacement.Text = \"Basic\"\n[void]$findNovice.Execute($findNovice.Text, $false, $false, $false, $false, $false, $true, 1, $false, $findNovice.Replacement.Text, 2)\n\n# 2. \"Beginner\" -> \"Basic\" (Linux bash, JavaScript, Ruby on Rails)\n$findBeginner = $d.Content.Find\n$findBeginner.Text = \"Beginner\"\n$findBeginner.Replacement.Text = \"Basic\"\n[void]$findBeginner.Execute($findBeginner.Text, $false, $false, $false, $false, $false, $true, 1, $false, $findBeginner.Replacement.Text, 2)\n\n# 3. Move the \"_GoBack\" bookmark (present at the top of the resume) to sit\n#    right after \"Ruby on Rails (Basic\" -- the location of the last text\n#    edit -- which is where Word leaves it after this kind of change.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$rng = $d.Content\n$findRuby = $rng.Find\n$findRuby.Text = \"Ruby on Rails (Basic\"\n$found = $findRuby.Execute()\nif ($found) {\n  $rng.Collapse(0) # wdCollapseEnd\n  [void]$d.Bookmarks.Add(\"_GoBack\", $rng)\n}\n"}
